# Update the default folder path values on the "Default_folders" sheet
# to reflect the author's new machine / folder layout.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Default_folders")

$ws.Range("B4").Value = "D:/Projects_data"
$ws.Range("B5").Value = "B:/Proteomics/Archive"
$ws.Range("B6").Value = "B:/Proteomics/Projects"
$ws.Range("B7").Value = "D:/Organisms"

# Leave the active selection where the author last left it when saving.
$ws.Range("B10").Select() | Out-Null
